{"js": "// The chart in the document is embedded as an inline picture. Its\n// description (alt text) doubles as a TinyButStrong/OpenTBS merge tag.\n// This edit tags the chart with the \"[chart_evol_by_cat]\" TBS field by\n// prefixing it onto the existing human-readable description, and clears\n// the (now redundant) alt-text title.\nconst inlinePictures = context.document.body.inlinePictures;\ninlinePictures.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < inlinePictures.items.length; i++) {\n  const pic = inlinePictures.items[i];\n  pic.load([\"altTextDescription\", \"altTextTitle\"]);\n}\nawait context.sync();\n\nfor (let i = 0; i < inlinePictures.items.length; i++) {\n  const pic = inlinePictures.items[i];\n  if (pic.altTextDescription === \"This is just a nice chart\") {\n    pic.altTextDescription = \"[chart_evol_by_cat]\\n\\nThis is just a nice chart\";\n    pic.altTextTitle = \"\";\n  }\n}\nawait context.sync();\n", "ps1": "# The chart in the document is embedded as an inline shape (InlineShapes\n# collection). Its AlternativeText (the \"Description\" field of the\n# picture's Alt Text pane) doubles as a TinyButStrong/OpenTBS merge tag.\n# Tag the chart with the \"[chart_evol_by_cat]\" TBS field by prefixing it\n# onto the existing human-readable description, and clear the (now\n# redundant) alt-text Title.\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le $d.InlineShapes.Count; $i++) {\n    $shp = $d.InlineShapes.Item($i)\n    if ($shp.AlternativeText -eq \"This is just a nice chart\") {\n        $shp.AlternativeText = \"[chart_evol_by_cat]`n`nThis is just a nice chart\"\n        $shp.Title = \"\"\n    }\n}\n"}
